$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''23.448.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '''1.644.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '''1.000'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '''300.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("D7").Value = '''0.3792'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("D8").Value = '''50.57'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("D9").Value = '''0.3497'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("D10").Value = '''0.08063'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").Value = '''1.214'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("D12").Value = '''1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("D13").Value = '''22.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").Value = '''6.299'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("D15").Value = '''7.238'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '''1.641.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = '''94.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").Value = '''0.06968'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").Value = '''6.615'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").Value = '''17.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").Value = '''12.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("D24").Value = '''23.466.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '''2.424'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.77%  '
$ws.Range("D26").Value = '''2.956'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").Value = '''20.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").Value = '''151.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("D29").Value = '''5.186'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").Value = '''131.70'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("D31").Value = '''1.825.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("D32").Value = '''6.821'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.25%  '
$ws.Range("D33").Value = '''2.139'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.78%  '
$ws.Range("D34").Value = '''11.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.26%  '
$ws.Range("D35").Value = '''0.9868'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.65%  '
$ws.Range("D36").Value = '''0.02681'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.90%  '
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").Value = '''5.898'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("D39").Value = '''0.2417'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.06%  '
$ws.Range("D40").Value = '''0.06776'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("D41").Value = '''12.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("D42").Value = '''0.6862'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").Value = '''1.293'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").Value = '''15.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("D45").Value = '''0.9995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").Value = '''0.6376'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '''3.924'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''2.242'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.07672'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.45%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''126.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("E51").Value = '  +2.22%  '
